$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.744734287261963
$ws.Range("B1").Value = 4.632020473480225
$ws.Range("C1").Value = 2.864572525024414
$ws.Range("D1").Value = 1.464893460273743
$ws.Range("E1").Value = 1.080158352851868
